# Insert a new data row at row 134 (pushes existing rows 134-251 down to
# 135-252) and populate it with the new record, matching the commit's
# "fruta / hortaliza, semanal" weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(134).Insert()

$ws.Range("A134").Value = 6
$ws.Range("B134").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C134").Value = "Metropolitana"
$ws.Range("D134").Value = 44729
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = 100112026
$ws.Range("G134").Value = "Haba"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 550
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 20000
$ws.Range("M134").Value = 19091
$ws.Range("N134").Value = "$/saco 25 kilos"
$ws.Range("O134").Value = "Provincia del Elquí"
$ws.Range("P134").Value = 764
$ws.Range("Q134").Value = 25
$ws.Range("R134").Value = "Hortaliza"
